# offresonline_tenders.xlsx update
# - Remove the tender "Centre hospitalier regional d'agadir" (accueil, long wording) due 19/06/2025
# - Remove the tender "Societe de developpement local rabat sale temara" (abattoirs bouknadel) due 23/06/2025
# - Add a new tender "Centre hospitalier provincial de kenitra" (accueil malades) due 07/07/2025,
#   placed right before the "Centre hospitalier provincial de tetouan" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like string into a cell as TEXT (not auto-converted to a date
# serial number) by routing it through a formula (string literal) + paste-values,
# using a scratch cell well outside the used range so it never shows up afterwards.
function Set-TextValue($cell, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# 1) Drop the old-wording Agadir "accueil" tender (row 2 / 19-06-2025).
$ws.Rows("2:2").Delete()

# 2) Drop the Rabat-Sale-Temara "abattoirs de bouknadel" tender.
#    After the first deletion it now sits on row 4.
$ws.Rows("4:4").Delete()

# 3) Insert the new Kenitra tender right before "Centre hospitalier provincial de
#    tetouan", which is now on row 8.
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "Organisme : Centre hospitalier provincial de kenitra. Objet : Appel d'offres ouvert a majoration : activit" + [char]0x00E9 + "s d" + [char]0x2019 + "accueil des malades du centre hospitalier provincial de kenitra."
Set-TextValue $ws.Range("B8") "07/07/2025"

Write-Host "Done"
